$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "contenuto" column (C) for the weeks that already took place
#    (rows 33-38). The previously "planned" text is replaced by what was
#    actually covered in class, the week 7 "current" row content moves down
#    to week 8, and the now-past highlighted session (row 38) gets its
#    exercise label filled in.
# ---------------------------------------------------------------------------
$ws.Range("C34").Value = "ripresa bi-dimensional gridding; rope e esercizio classifier; presentazione assignment da Marco"
$ws.Range("C35").Value = "normal normal in pymc3  (esercizio su predictive;  spiegato fino a student likelihood, inclusa)"
$ws.Range("C36").Value = "exe MF: metropolis"
$ws.Range("C37").Value = "finire normal normal model; presentazione progetto; intro hyp test"
$ws.Range("C38").Value = "exe"
$ws.Range("C33").Value = "?"

# ---------------------------------------------------------------------------
# 2. Move the "current session" highlight down by one row: it used to mark
#    row 35, it now marks row 36. Copying the formatting keeps number
#    formats / fills consistent with the rest of the table, then row 35 is
#    restored to the plain (non-highlighted) look used by the surrounding
#    rows.
# ---------------------------------------------------------------------------
$ws.Range("A35:C35").Copy() | Out-Null
$ws.Range("A36:C36").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("A34:C34").Copy() | Out-Null
$ws.Range("A35:C35").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Refresh the view: scrolled up one row and the active selection moved to
#    the next editable cell below the table.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
$ws.Range("C39").Select()

# Widen the workbook window, matching the wider layout used afterwards.
$win.Width = 32000
